# Rename the embedded logo pictures' display names:
#   - the two Pearson Edexcel logo pictures (in the primary and
#     first-page footers) from "image1.png" to "image2.png"
#   - the BTEC logo picture (in the first-page header) from
#     "image2.jpg" to "image1.jpg"
#
# InlineShape objects don't expose a writable Name in Word's object
# model, so each picture is briefly converted to a floating Shape
# (which does expose Name), renamed, then converted back to an
# inline shape in place.

$d = $word.ActiveDocument

function Rename-InlinePicture($range, $newName) {
    for ($j = 1; $j -le $range.InlineShapes.Count; $j++) {
        $shp = $range.InlineShapes.Item($j).ConvertToShape()
        $shp.Name = $newName
        $shp.ConvertToInlineShape() | Out-Null
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    # Primary footer (maps to the picture with id="2") and first-page
    # footer (id="3") both carry the Pearson Edexcel logo.
    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            Rename-InlinePicture $ftr.Range "image2.png"
        }
    }

    # First-page header (id="1") carries the BTEC logo.
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            Rename-InlinePicture $hdr.Range "image1.jpg"
        }
    }
}
